$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.281.87'
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").Value = '3.258.62'
$ws.Range("E3").Value = '  +5.37%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '578.00'
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").Value = '154.21'
$ws.Range("E6").Value = '  +6.56%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.249.91'
$ws.Range("E8").Value = '  +5.53%  '
$ws.Range("E9").Value = '  +3.79%  '
$ws.Range("E10").Value = '  +7.43%  '
$ws.Range("E11").Value = '  +3.88%  '
$ws.Range("D12").Value = '0.493'
$ws.Range("E12").Value = '  +3.44%  '
$ws.Range("D13").Value = '37.91'
$ws.Range("E13").Value = '  +3.13%  '
$ws.Range("E14").Value = '  +4.78%  '
$ws.Range("D15").Value = '3.780.69'
$ws.Range("E15").Value = '  +5.25%  '
$ws.Range("D16").Value = '562.76'
$ws.Range("E16").Value = '  +11.49%  '
$ws.Range("D17").Value = '66.316.60'
$ws.Range("E17").Value = '  +2.33%  '
$ws.Range("D18").Value = '3.259.52'
$ws.Range("E18").Value = '  +5.45%  '
$ws.Range("E19").Value = '  +2.60%  '
$ws.Range("D20").Value = '7.15'
$ws.Range("E20").Value = '  +5.04%  '
$ws.Range("D21").Value = '14.50'
$ws.Range("E21").Value = '  +3.67%  '
$ws.Range("D22").Value = '0.746'
$ws.Range("E22").Value = '  +6.60%  '
$ws.Range("D23").Value = '7.85'
$ws.Range("E23").Value = '  +7.19%  '
$ws.Range("D24").Value = '13.64'
$ws.Range("E24").Value = '  +5.24%  '
$ws.Range("D25").Value = '82.50'
$ws.Range("E25").Value = '  +3.00%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '9.38'
$ws.Range("E27").Value = '  +17.57%  '
$ws.Range("E29").Value = '  +5.11%  '
$ws.Range("D30").Value = '27.99'
$ws.Range("E30").Value = '  +5.03%  '
$ws.Range("D31").Value = '2.76'
$ws.Range("E31").Value = '  +2.16%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").Value = '  +3.86%  '
$ws.Range("D34").Value = '569.18'
$ws.Range("E34").Value = '  +9.79%  '
$ws.Range("D35").Value = '5.77'
$ws.Range("E35").Value = '  +3.66%  '
$ws.Range("D36").Value = '6.43'
$ws.Range("E36").Value = '  +4.68%  '
$ws.Range("D37").Value = '0.0461'
$ws.Range("E37").Value = '  +12.02%  '
$ws.Range("D38").Value = '55.23'
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("E39").Value = '  +6.51%  '
$ws.Range("D40").Value = '3.09'
$ws.Range("E40").Value = '  +12.94%  '
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").Value = '3.165.90'
$ws.Range("E42").Value = '  +6.52%  '
$ws.Range("D43").Value = '8.65'
$ws.Range("E43").Value = '  +1.66%  '
$ws.Range("D44").Value = '0.276'
$ws.Range("E44").Value = '  +9.80%  '
$ws.Range("D45").Value = '2.30'
$ws.Range("E45").Value = '  +5.34%  '
$ws.Range("D46").Value = '26.67'
$ws.Range("E46").Value = '  +3.68%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").Value = '0.0₃0558'
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("D49").Value = '125.06'
$ws.Range("E49").Value = '  +3.51%  '
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("E51").Value = '  +6.99%  '
